$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "password123" entry for each non-admin user (rows 3-6, column B)
$ws.Range("B3").Value = "password123"
$ws.Range("B4").Value = "password123"
$ws.Range("B5").Value = "password123"
$ws.Range("B6").Value = "password123"

# Paolo Cisneros drops from access level 3 to 2
$ws.Range("C5").Value = 2

# Widen column B to fit the new password values (OOXML width ends up 5/6
# wider than the ColumnWidth value we set, so back that offset out here)
$ws.Range("B:B").ColumnWidth = 11.666666666666666

# Leave selection on A5 (matches author's last click before saving)
$ws.Range("A5").Select()
